# Update "想去人数" (F column) counts for several camp/show events across
# the workbook's sheets, matching the refreshed data pull from bilibili.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3644
$ws1.Range("F5").Value  = 3644
$ws1.Range("F7").Value  = 5171
$ws1.Range("F22").Value = 4950
$ws1.Range("F26").Value = 6076
$ws1.Range("F32").Value = 4449
$ws1.Range("F36").Value = 1061
$ws1.Range("F40").Value = 886

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 28

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1129

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1129
$ws4.Range("F7").Value  = 3644
$ws4.Range("F8").Value  = 3644
$ws4.Range("F10").Value = 5171
$ws4.Range("F26").Value = 4950
$ws4.Range("F30").Value = 6076
$ws4.Range("F36").Value = 4449
$ws4.Range("F38").Value = 28
$ws4.Range("F41").Value = 1061
$ws4.Range("F45").Value = 886
